# Updates cryptos list values (prices/volume %) per scraper refresh.
# For Price (column D) values that look like plain numbers (e.g. "1.001"),
# a leading apostrophe is used so Excel stores them as text -- matching the
# original inlineStr/text cell type in the workbook -- instead of silently
# converting them to floating point numbers (which would also lose trailing
# zeros / precision, e.g. "1.0000" -> 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.159.93"
$ws.Range("E2").Value = "  -0.59%  "

# Row 3
$ws.Range("D3").Value = "1.829.06"
$ws.Range("E3").Value = "  +0.92%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").Value = "'310.94"
$ws.Range("E5").Value = "  -0.71%  "

# Row 6
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.16%  "

# Row 7
$ws.Range("D7").Value = "'0.4966"
$ws.Range("E7").Value = "  -3.77%  "

# Row 8
$ws.Range("D8").Value = "'0.3928"
$ws.Range("E8").Value = "  -2.64%  "

# Row 9
$ws.Range("D9").Value = "'0.09849"
$ws.Range("E9").Value = "  +25.06%  "

# Row 10
$ws.Range("D10").Value = "'1.112"
$ws.Range("E10").Value = "  -0.27%  "

# Row 11
$ws.Range("D11").Value = "'41.07"
$ws.Range("E11").Value = "  -0.16%  "

# Row 12
$ws.Range("D12").Value = "'6.474"
$ws.Range("E12").Value = "  +1.60%  "

# Row 13
$ws.Range("D13").Value = "'20.67"
$ws.Range("E13").Value = "  +0.89%  "

# Row 14
$ws.Range("E14").Value = "  +0.17%  "

# Row 15
$ws.Range("D15").Value = "1.818.32"
$ws.Range("E15").Value = "  +1.42%  "

# Row 16
$ws.Range("D16").Value = "'7.313"
$ws.Range("E16").Value = "  -0.58%  "

# Row 17
$ws.Range("D17").Value = "'0.00001144"
$ws.Range("E17").Value = "  +5.58%  "

# Row 18
$ws.Range("D18").Value = "'92.72"
$ws.Range("E18").Value = "  -0.19%  "

# Row 19
$ws.Range("D19").Value = "'0.06657"
$ws.Range("E19").Value = "  +0.94%  "

# Row 20
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.20%  "

# Row 21
$ws.Range("D21").Value = "'17.27"
$ws.Range("E21").Value = "  -0.59%  "

# Row 22
$ws.Range("D22").Value = "'5.994"
$ws.Range("E22").Value = "  -0.91%  "

# Row 23
$ws.Range("D23").Value = "28.199.95"
$ws.Range("E23").Value = "  -0.65%  "

# Row 24
$ws.Range("D24").Value = "'11.36"
$ws.Range("E24").Value = "  +1.28%  "

# Row 25
$ws.Range("D25").Value = "'2.247"
$ws.Range("E25").Value = "  +0.94%  "

# Row 26
$ws.Range("D26").Value = "'158.87"
$ws.Range("E26").Value = "  -1.22%  "

# Row 27
$ws.Range("D27").Value = "'20.85"
$ws.Range("E27").Value = "  +0.98%  "

# Row 28
$ws.Range("D28").Value = "2.037.43"
$ws.Range("E28").Value = "  +1.11%  "

# Row 29
$ws.Range("D29").Value = "'2.430"
$ws.Range("E29").Value = "  +0.73%  "

# Row 30
$ws.Range("D30").Value = "'127.05"
$ws.Range("E30").Value = "  -1.39%  "

# Row 31
$ws.Range("D31").Value = "'0.1059"
$ws.Range("E31").Value = "  -2.56%  "

# Row 32
$ws.Range("D32").Value = "'1.041"
$ws.Range("E32").Value = "  -1.30%  "

# Row 33
$ws.Range("D33").Value = "'5.612"
$ws.Range("E33").Value = "  +0.23%  "

# Row 34
$ws.Range("D34").Value = "'3.618"
$ws.Range("E34").Value = "  -1.10%  "

# Row 35
$ws.Range("D35").Value = "'0.06744"
$ws.Range("E35").Value = "  -6.33%  "

# Row 36
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'9.001"
$ws.Range("E36").Value = "  -1.73%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02353"
$ws.Range("E37").Value = "  +0.74%  "

# Row 38
$ws.Range("D38").Value = "'0.2151"
$ws.Range("E38").Value = "  -0.78%  "

# Row 39
$ws.Range("E39").Value = "  -1.88%  "

# Row 40
$ws.Range("D40").Value = "'4.985"
$ws.Range("E40").Value = "  -1.81%  "

# Row 41
$ws.Range("D41").Value = "'0.6237"
$ws.Range("E41").Value = "  -0.01%  "

# Row 42
$ws.Range("E42").Value = "  +1.83%  "

# Row 43
$ws.Range("D43").Value = "'1.0000"
$ws.Range("E43").Value = "  +0.15%  "

# Row 44
$ws.Range("D44").Value = "'13.24"
$ws.Range("E44").Value = "  -0.40%  "

# Row 45
$ws.Range("D45").Value = "'0.5950"
$ws.Range("E45").Value = "  -1.21%  "

# Row 46
$ws.Range("D46").Value = "'3.704"
$ws.Range("E46").Value = "  -1.14%  "

# Row 47
$ws.Range("D47").Value = "'1.276"
$ws.Range("E47").Value = "  -3.56%  "

# Row 48
$ws.Range("D48").Value = "'124.33"
$ws.Range("E48").Value = "  -1.55%  "

# Row 49
$ws.Range("D49").Value = "'1.953"
$ws.Range("E49").Value = "  +0.43%  "

# Row 50
$ws.Range("D50").Value = "'1.183"
$ws.Range("E50").Value = "  -2.81%  "

# Row 51
$ws.Range("D51").Value = "'0.06790"
$ws.Range("E51").Value = "  -1.04%  "

